$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric/percent-looking cells so Excel keeps them as literal strings
$textCells = @("D2", "E2", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "D22", "E22", "E23", "D24", "E24", "E25", "D26", "E26", "D38", "D39", "E39", "D40", "E40", "D41", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47", "D48", "E48", "E49", "D50", "E50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "292.13"
$ws.Range("E2").Value = "-2.68%"
$ws.Range("E3").Value = "-0.84%"
$ws.Range("D4").Value = "5.013"
$ws.Range("E4").Value = "-2.63%"
$ws.Range("D5").Value = "0.07291"
$ws.Range("E5").Value = "-2.82%"
$ws.Range("D6").Value = "1.528"
$ws.Range("E6").Value = "-6.43%"
$ws.Range("D7").Value = "0.9262"
$ws.Range("E7").Value = "-1.47%"
$ws.Range("D8").Value = "2.351"
$ws.Range("E8").Value = "-3.01%"
$ws.Range("D9").Value = "0.1155"
$ws.Range("E9").Value = "-3.66%"
$ws.Range("D10").Value = "0.1766"
$ws.Range("E10").Value = "-0.87%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.08787"
$ws.Range("E11").Value = "-1.52%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.04353"
$ws.Range("E12").Value = "4.99%"
$ws.Range("D13").Value = "0.1053"
$ws.Range("E13").Value = "-0.07%"
$ws.Range("D14").Value = "0.001272"
$ws.Range("E14").Value = "-1.51%"
$ws.Range("D15").Value = "0.005975"
$ws.Range("E15").Value = "2.90%"
$ws.Range("E16").Value = "-0.22%"
$ws.Range("D17").Value = "4.285"
$ws.Range("E17").Value = "-1.33%"
$ws.Range("D18").Value = "0.3280"
$ws.Range("E18").Value = "-2.17%"
$ws.Range("D19").Value = "7.955"
$ws.Range("E19").Value = "4.18%"
$ws.Range("D20").Value = "0.1391"
$ws.Range("E20").Value = "2.61%"
$ws.Range("E21").Value = "-1.44%"
$ws.Range("D22").Value = "0.03925"
$ws.Range("E22").Value = "1.81%"
$ws.Range("E23").Value = "-1.92%"
$ws.Range("D24").Value = "0.003653"
$ws.Range("E24").Value = "-7.82%"
$ws.Range("E25").Value = "-7.86%"
$ws.Range("D26").Value = "0.0003728"
$ws.Range("E26").Value = "-0.23%"
$ws.Range("D38").Value = "0.02303"
$ws.Range("D39").Value = "0.05044"
$ws.Range("E39").Value = "-0.51%"
$ws.Range("D40").Value = "0.005860"
$ws.Range("E40").Value = "67.10%"
$ws.Range("D41").Value = "0.007846"
$ws.Range("E41").Value = "1.39%"
$ws.Range("E42").Value = "-0.85%"
$ws.Range("D43").Value = "0.007375"
$ws.Range("E43").Value = "-2.73%"
$ws.Range("D44").Value = "0.007234"
$ws.Range("E44").Value = "-9.60%"
$ws.Range("D45").Value = "0.3189"
$ws.Range("E45").Value = "-1.96%"
$ws.Range("D46").Value = "0.00006214"
$ws.Range("E46").Value = "-8.94%"
$ws.Range("E47").Value = "-0.23%"
$ws.Range("D48").Value = "0.03890"
$ws.Range("E48").Value = "-84.49%"
$ws.Range("E49").Value = "-0.23%"
$ws.Range("D50").Value = "0.0002003"
$ws.Range("E50").Value = "-0.23%"
